# Weekly driver report update for 2025-04-29
# Updates the "Good Drivers" table on the Driver Summary sheet: new headers
# (columns A-J) and a refreshed data row, plus matching column widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths ---------------------------------------------------
# ColumnWidth (characters) maps to stored OOXML width as N + 5/6, so back
# the offset out to land on the exact target widths from the diff.
function Set-ExactColumnWidth($col, $target) {
    $ws.Columns.Item($col).ColumnWidth = $target - (5/6)
}

Set-ExactColumnWidth 2 14    # B: 15 -> 14
Set-ExactColumnWidth 5 14    # E: 16 -> 14
Set-ExactColumnWidth 6 11    # F: 2  -> 11
Set-ExactColumnWidth 7 48    # G: 2  -> 48
Set-ExactColumnWidth 8 15    # H: 2  -> 15
Set-ExactColumnWidth 9 30    # I: 2  -> 30
Set-ExactColumnWidth 10 16   # J: 2  -> 16

# --- Row 11: replace the "Good Drivers" header row --------------------
# Old header row had bold/bordered styling (s="2"/s="3"); new header row
# is plain, unstyled text across columns A-J.
$headerRow = $ws.Range("A11:J11")
$headerRow.ClearFormats()

$ws.Range("A11").Value = "adapter-driver"
$ws.Range("B11").Value = "good sum"
$ws.Range("C11").Value = "critical sum"
$ws.Range("D11").Value = "warning sum"
$ws.Range("E11").Value = "client count"
$ws.Range("F11").Value = "total sum"
$ws.Range("G11").Value = "adapter"
$ws.Range("H11").Value = "driver"
$ws.Range("I11").Value = "good roaming calculation (%)"
$ws.Range("J11").Value = "driver vintage"

# --- Row 12: new data row (previously blank) ---------------------------
$ws.Range("A12").Value = "Realtek RTL8852AE WiFi 6 802.11ax PCIe Adapter - 6001.10.356.1"
$ws.Range("B12").Value = 1071383
$ws.Range("C12").Value = 4419
$ws.Range("D12").Value = 180
$ws.Range("E12").Value = 1644
$ws.Range("F12").Value = 1075982
$ws.Range("G12").Value = "realtek rtl8852ae wifi 6 802.11ax pcie adapter"
$ws.Range("H12").Value = "6001.10.356.1"
$ws.Range("I12").Value = 99.59999999999999

# J12 looks like a date ("2024-05-12"); force text so Excel doesn't
# auto-convert it to a date serial number, then strip the number-format
# artifact that forcing text leaves behind so the cell stays unstyled.
$ws.Range("J12").NumberFormat = "@"
$ws.Range("J12").Value = "2024-05-12"
$ws.Range("J12").ClearFormats()

# --- Rows 13-16 are no longer part of the used range --------------------
$ws.Range("A13:J16").Clear()
